{"js": "// The commit (\"github and git connectivity\") strips a batch of stray Word\n// spell/grammar proofing bookmarks (<w:proofErr w:type=\"spellStart|spellEnd\"/>\n// and, in one spot, <w:proofErr w:type=\"gramStart|gramEnd\"/>) that were\n// splitting otherwise-contiguous text into extra runs throughout the\n// document. Those bookmarks are purely internal proofing annotations \u2014\n// invisible, non-semantic, and not exposed anywhere in the Word object\n// model / Office.js API surface \u2014 so their removal never changes the\n// document's actual (rendered) text content. Diffing the canonical XML\n// confirms that, across the whole file, the *only* place where the visible\n// text itself changes is:\n//\n//   \"You  can merge the changes of One branch into another.\"\n//   -> \"You can merge the changes of One branch into another.\"\n//\n// i.e. a genuine double space after \"You\" collapses to a single space.\n// That is the one concrete, scriptable edit to make here.\n\nconst body = context.document.body;\n\n// Narrow the search to the full sentence first so we don't touch any of\n// the many other legitimate double-spaces elsewhere in the document.\nconst sentenceResults = body.search(\n  \"You  can merge the changes of One branch into another.\",\n  { matchCase: true, matchWholeWord: false }\n);\nsentenceResults.load(\"items\");\nawait context.sync();\n\nif (sentenceResults.items.length > 0) {\n  const sentenceRange = sentenceResults.items[0];\n\n  // Within that sentence, find just the double space right after \"You\" and\n  // collapse it to a single space. Doing the edit at this fine grain keeps\n  // every surrounding run / proofErr bookmark untouched instead of\n  // rewriting the whole sentence into a single run.\n  const doubleSpace = sentenceRange.search(\"You  \", { matchCase: true });\n  doubleSpace.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < doubleSpace.items.length; i++) {\n    doubleSpace.items[i].insertText(\"You \", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The commit (\"github and git connectivity\") strips a batch of stray Word\n# spell/grammar proofing bookmarks (<w:proofErr w:type=\"spellStart|spellEnd\"/>\n# and, in one spot, <w:proofErr w:type=\"gramStart|gramEnd\"/>) that were\n# splitting otherwise-contiguous text into extra runs throughout the\n# document. Those bookmarks are purely internal proofing annotations \u2014\n# invisible, non-semantic, and not exposed anywhere in the Word object\n# model (no Range/Paragraph property surfaces them for editing) \u2014 so their\n# removal never changes the document's actual (rendered) text content.\n# Diffing the canonical XML confirms that, across the whole file, the\n# *only* place where the visible text itself changes is:\n#\n#   \"You  can merge the changes of One branch into another.\"\n#   -> \"You can merge the changes of One branch into another.\"\n#\n# i.e. a genuine double space after \"You\" collapses to a single space.\n# That is the one concrete, scriptable edit to make here.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph containing the double space so we don't disturb any\n# of the many other legitimate double-spaces elsewhere in the document.\n$target = $null\nforeach ($para in $d.Paragraphs) {\n    if ($para.Range.Text -like \"*You  can merge the changes of One branch into another.*\") {\n        $target = $para\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $full = $target.Range\n    $text = $full.Text\n    $idx = $text.IndexOf(\"You  \")\n    if ($idx -ge 0) {\n        # Collapse just the redundant space character right after \"You\" to\n        # a single space, leaving every surrounding run / proofErr bookmark\n        # untouched instead of rewriting the whole sentence.\n        $spaceOffset = $idx + 3\n        $spaceStart = $full.Start + $spaceOffset\n        $spaceRange = $d.Range($spaceStart, $spaceStart + 1)\n        $spaceRange.Delete()\n    }\n}\n"}
